# Fruta / hortaliza, semanal
# Apply a cyclic shift of weekly data among rows 4-7:
#   row4 -> row7, row5 -> row4, row6 -> row5, row7 -> row6
# Only columns D, L, M, N, O, P, S differ between rows; update those values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original (pre-edit) values for rows 4-7 in the affected columns.
$cols = @("D", "L", "M", "N", "O", "P", "S")

$orig = @{}
foreach ($row in 4..7) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value()
    }
}

# New row order: target row gets the values that were previously in the next
# row, with row 7 wrapping around to take row 4's original values.
$mapping = @{
    4 = 5
    5 = 6
    6 = 7
    7 = 4
}

foreach ($targetRow in 4..7) {
    $sourceRow = $mapping[$targetRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $orig[$sourceRow][$col]
    }
}
